# Applies the changes described by the diff to the "Export" sheet:
#   - account 004526450 (MSD) balance changes 30000 -> 62000 and moves up
#     (sheet is sorted by balance, descending) to just above account 004213929
#   - account 004479287 (ANA) balance changes 30642.81 -> 32392.83 and moves up
#     to just above account 004451978
#   - three brand new rows are added:
#       005046919 MARIANA 3740   (just above account 005000460)
#       004504449 KELMA    800   (just above account 004462930)
#       004979322 MARILIA  500   (just below account 004405234 / JULIO)
#
# Because row numbers shift whenever a row is inserted or deleted, every
# operation below is carried out from the bottom of the sheet upward so
# that none of the row numbers determined ahead of time are invalidated
# by an earlier step.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Small helper: make sure we're about to insert/delete next to the row we
# think we are, before mutating anything (row numbers are hard-coded below,
# derived from the sheet's current, known layout).
function Assert-Account($row, $expected) {
  $actual = $ws.Cells.Item($row, 1).Value2
  if ($actual -ne $expected) {
    throw "row $row : expected account $expected but found $actual"
  }
}

# Account numbers are stored as text with significant leading zeros
# ("004526450", not 4526450). Typing a leading-zero value straight into a
# General-formatted cell makes Excel reinterpret it as a number, so every
# account number below is written with a leading apostrophe - the normal
# Excel trick for forcing text entry - which preserves the digits exactly.

# 1) New row: 004979322 / MARILIA / 500, inserted directly below account
#    004405234 (JULIO, 522.09) i.e. directly above account 004547722 (MARCIA).
Assert-Account 177 "004405234"
Assert-Account 178 "004547722"
$ws.Rows.Item(178).Insert()
$ws.Cells.Item(178,1).Value = "'004979322"
$ws.Cells.Item(178,2).Value = "MARILIA"
$ws.Cells.Item(178,3).Value = 500

# 2) New row: 004504449 / KELMA / 800, inserted directly above account 004462930.
Assert-Account 174 "004462930"
$ws.Rows.Item(174).Insert()
$ws.Cells.Item(174,1).Value = "'004504449"
$ws.Cells.Item(174,2).Value = "KELMA"
$ws.Cells.Item(174,3).Value = 800

# 3) New row: 005046919 / MARIANA / 3740, inserted directly above account 005000460.
Assert-Account 168 "005000460"
$ws.Rows.Item(168).Insert()
$ws.Cells.Item(168,1).Value = "'005046919"
$ws.Cells.Item(168,2).Value = "MARIANA"
$ws.Cells.Item(168,3).Value = 3740

# 4) Remove the old 004526450 / MSD / 30000 row (it reappears, with a new
#    balance, further up the sheet - see step 6).
Assert-Account 99 "004526450"
$ws.Rows.Item(99).Delete()

# 5) Remove the old 004479287 / ANA / 30642.81 row (it reappears, with a new
#    balance, further up the sheet - see step 7).
Assert-Account 77 "004479287"
$ws.Rows.Item(77).Delete()

# 6) Re-insert 004479287 / ANA with its updated balance (32392.83), directly
#    above account 004451978.
Assert-Account 74 "004451978"
$ws.Rows.Item(74).Insert()
$ws.Cells.Item(74,1).Value = "'004479287"
$ws.Cells.Item(74,2).Value = "ANA"
$ws.Cells.Item(74,3).Value = 32392.83

# 7) Re-insert 004526450 / MSD with its updated balance (62000), directly
#    above account 004213929.
Assert-Account 35 "004213929"
$ws.Rows.Item(35).Insert()
$ws.Cells.Item(35,1).Value = "'004526450"
$ws.Cells.Item(35,2).Value = "MSD"
$ws.Cells.Item(35,3).Value = 62000
